# Exploit for case 12 and 28
# Fill in the newly-discovered testcase numbers (column E) across the
# bug-tracking sheet, and record the attack-script name / flag / student id
# for case12.py (row 7) and case28.py (row 15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Testcase numbers that were filled in (column E) ---
$ws.Range("E4").Value = 19
$ws.Range("E7").Value = 12
$ws.Range("E15").Value = 28
$ws.Range("E16").Value = 30
$ws.Range("E22").Value = 16
$ws.Range("E23").Value = 26
$ws.Range("E24").Value = 9
$ws.Range("E30").Value = 17
$ws.Range("E31").Value = 20
$ws.Range("E32").Value = 23

# --- case12.py details (row 7) ---
# Set the flag column first so the shared-string table picks up the new
# unique strings in the same order the original author's workbook did.
$ws.Range("G7").Value = "ccCL2uX5L4kGU52"
$ws.Range("F7").Value = "case12.py"
$ws.Range("H7").Value = "A0127604L"

# --- case28.py details (row 15) ---
$ws.Range("F15").Value = "case28.py"
$ws.Range("G15").Value = "QG3PwQjJmsNnQrx"
$ws.Range("H15").Value = "A0127604L"

# --- Update the view so the newly-added rows are visible/selected ---
$ws.Range("G31").Select()
